$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any existing rows beyond the new data range (old data went to row 81, new ends at row 64)
$ws.Range("A65:E81").ClearContents()

# Target data rows (Dia, total_venda, Mes, Ano, Periodo)
$data = @(
    @(28,8693.200000000001,7,2025,"07/2025"),
    @(2,45067.52,6,2025,"06/2025"),
    @(3,35136.72,6,2025,"06/2025"),
    @(4,11872.2,6,2025,"06/2025"),
    @(5,12749.74,6,2025,"06/2025"),
    @(6,12978.82,6,2025,"06/2025"),
    @(9,11059.45,6,2025,"06/2025"),
    @(10,4973.1,6,2025,"06/2025"),
    @(11,19810.51,6,2025,"06/2025"),
    @(12,17504.38,6,2025,"06/2025"),
    @(13,9230.379999999999,6,2025,"06/2025"),
    @(16,24913.77,6,2025,"06/2025"),
    @(17,14944.85,6,2025,"06/2025"),
    @(18,14248.93,6,2025,"06/2025"),
    @(20,5940.05,6,2025,"06/2025"),
    @(21,119.13,6,2025,"06/2025"),
    @(23,44131.86,6,2025,"06/2025"),
    @(24,15341.74,6,2025,"06/2025"),
    @(25,14231.97,6,2025,"06/2025"),
    @(26,24113.52,6,2025,"06/2025"),
    @(27,9113.01,6,2025,"06/2025"),
    @(30,81604.48,6,2025,"06/2025"),
    @(2,20463.44,5,2025,"05/2025"),
    @(5,29720.49,5,2025,"05/2025"),
    @(6,19666.55,5,2025,"05/2025"),
    @(7,21883.41,5,2025,"05/2025"),
    @(8,45418.95,5,2025,"05/2025"),
    @(9,39575.07,5,2025,"05/2025"),
    @(12,15499.7,5,2025,"05/2025"),
    @(13,10155.95,5,2025,"05/2025"),
    @(14,34413.69,5,2025,"05/2025"),
    @(15,28325.85,5,2025,"05/2025"),
    @(16,12000.74,5,2025,"05/2025"),
    @(19,17756.75,5,2025,"05/2025"),
    @(20,23593.88,5,2025,"05/2025"),
    @(21,9475.469999999999,5,2025,"05/2025"),
    @(22,27766.95,5,2025,"05/2025"),
    @(23,11639.15,5,2025,"05/2025"),
    @(26,28728.43,5,2025,"05/2025"),
    @(27,15807.63,5,2025,"05/2025"),
    @(28,24406.06,5,2025,"05/2025"),
    @(29,30466.72,5,2025,"05/2025"),
    @(30,23720.49,5,2025,"05/2025"),
    @(1,31020.03,4,2025,"04/2025"),
    @(2,48690.41,4,2025,"04/2025"),
    @(3,35050.25,4,2025,"04/2025"),
    @(4,7210.62,4,2025,"04/2025"),
    @(7,10042.97,4,2025,"04/2025"),
    @(8,7630.52,4,2025,"04/2025"),
    @(9,19095.69,4,2025,"04/2025"),
    @(10,16812.5,4,2025,"04/2025"),
    @(11,6701.92,4,2025,"04/2025"),
    @(14,38115.7,4,2025,"04/2025"),
    @(15,12545.33,4,2025,"04/2025"),
    @(16,14981.31,4,2025,"04/2025"),
    @(17,14359.25,4,2025,"04/2025"),
    @(22,18901.75,4,2025,"04/2025"),
    @(23,22062.13,4,2025,"04/2025"),
    @(24,17802.7,4,2025,"04/2025"),
    @(25,6150.450000000001,4,2025,"04/2025"),
    @(28,31344.88,4,2025,"04/2025"),
    @(29,33546.1,4,2025,"04/2025"),
    @(30,32403.27,4,2025,"04/2025")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}

Write-Output "rows written: $($row - 2)"
